$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update last_edited_time for rows 7-12 (shared string used across these rows)
$ws.Range("D7:D12").Value = "2024-07-31T18:24:00.000Z"

# Update numeric figures on row 7
$ws.Range("T7").Value = 64300000
$ws.Range("W7").Value = 332530000
$ws.Range("AA7").Value = 376038000
$ws.Range("AE7").Value = 708568000
$ws.Range("AH7").Value = 615568000
$ws.Range("AK7").Value = 93
$ws.Range("AN7").Value = 93000000
$ws.Range("AQ7").Value = 679868000
